$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.295.44'
$ws.Range('E2').Value = '  +6.38%  '
$ws.Range('D3').Value = '4.061.93'
$ws.Range('E3').Value = '  +6.61%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '525.33'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.97'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.711'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +18.89%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.766'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +9.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.178'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000337'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '49.35'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +20.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.92'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.10%  '
$ws.Range('D14').Value = '4.714.15'
$ws.Range('E14').Value = '  +6.40%  '
$ws.Range('D15').Value = '4.090.88'
$ws.Range('E15').Value = '  +6.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.34'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '21.07'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.51%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.23'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.04%  '
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').Value = '72.236.69'
$ws.Range('E20').Value = '  +6.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '440.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.94%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '101.86'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +18.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.62'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '15.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.15'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.46'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.02'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.74'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +8.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.87'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.37'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +20.83%  '
$ws.Range('E31').Value = '  +5.31%  '
$ws.Range('E32').Value = '  +7.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '682.01'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.78'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +14.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '67.44'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '42.97'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.72%  '
$ws.Range('D37').Value = '0.0₃0873'
$ws.Range('E37').Value = '  +7.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.433'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.157'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.46'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0509'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +8.82%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('E43').Value = '  -0.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.16'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.36%  '
$ws.Range('E45').Value = '  +13.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.75'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.45%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.45'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.62'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +14.77%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.13'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +8.18%  '
$ws.Range('B50').Value = 'FLOKI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000280'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +11.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.38'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.49%  '
